$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.458.31"
$ws.Range("D3").Value = "2.163.33"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'228.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'63.65"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'16.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "2.483.28"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "'22.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "'0.813"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "2.173.27"
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").Value = "39.456.83"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "'6.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").Value = "'71.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "0.0₃0849"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").Value = "'229.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'2.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").Value = "'172.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "'9.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").Value = "'19.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("D30").Value = "'1.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "'2.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.39%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").Value = "'7.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "'0.0621"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").Value = "'2.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'103.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "'17.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").Value = "1.523.58"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "'1.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.59%  "
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "'2.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0926"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").Value = "2.367.35"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("E51").Value = "  -0.68%  "
